# Correcciones en reglas del documento stock actual
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("L3").Value = 0

# --- Row 13 (becomes hidden) ---
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("P13").Value = 24
$ws.Range("Q13").Value = 0
$ws.Range("U13").Value = 0
$ws.Rows.Item(13).Hidden = $true

# --- Row 15 (becomes hidden) ---
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("P15").Value = 37
$ws.Range("Q15").Value = 0
$ws.Range("U15").Value = 0
$ws.Rows.Item(15).Hidden = $true

# --- Row 16 (was visible, becomes hidden; swaps talla data with row 17) ---
$ws.Range("C16").Value = "0I8LA30  "
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 10.43
$ws.Range("H16").Value = 4.17
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("P16").Value = 52
$ws.Range("Q16").Value = 0
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 0
$ws.Rows.Item(16).Hidden = $true

# --- Row 17 (remains hidden; swaps talla data with row 16) ---
$ws.Range("C17").Value = "5LA75    "
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 37.62
$ws.Range("H17").Value = 15.05
$ws.Range("L17").Value = 0
$ws.Range("P17").Value = 21
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 0

# --- Row 18 (becomes hidden) ---
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("P18").Value = 10
$ws.Range("Q18").Value = 0
$ws.Range("U18").Value = 0
$ws.Rows.Item(18).Hidden = $true

# --- Row 21 (stays hidden) ---
$ws.Range("L21").Value = 0

# --- Row 25 (stays visible) ---
$ws.Range("M25").Value = 94.88
$ws.Range("N25").Value = 56.93
$ws.Range("P25").Value = 2
$ws.Range("Q25").Value = 5
$ws.Range("U25").Value = 5

# --- Row 26 (becomes hidden) ---
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("P26").Value = 12
$ws.Range("Q26").Value = 0
$ws.Range("U26").Value = 0
$ws.Rows.Item(26).Hidden = $true

# --- Row 27 (becomes hidden) ---
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("P27").Value = 30
$ws.Range("Q27").Value = 0
$ws.Range("U27").Value = 0
$ws.Rows.Item(27).Hidden = $true

# --- Row 28 (becomes hidden) ---
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("P28").Value = 9
$ws.Range("Q28").Value = 0
$ws.Range("U28").Value = 0
$ws.Rows.Item(28).Hidden = $true

# --- Row 29 (becomes hidden) ---
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("P29").Value = 19
$ws.Range("Q29").Value = 0
$ws.Range("U29").Value = 0
$ws.Rows.Item(29).Hidden = $true

# --- Row 30 (becomes hidden) ---
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("P30").Value = 28
$ws.Range("Q30").Value = 0
$ws.Range("U30").Value = 0
$ws.Rows.Item(30).Hidden = $true

# --- Row 31 (becomes hidden) ---
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("P31").Value = 15
$ws.Range("Q31").Value = 0
$ws.Range("U31").Value = 0
$ws.Rows.Item(31).Hidden = $true

# --- Row 34 (becomes hidden) ---
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("P34").Value = 4
$ws.Range("Q34").Value = 0
$ws.Range("U34").Value = 0
$ws.Rows.Item(34).Hidden = $true

# --- Row 35 (becomes hidden) ---
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("P35").Value = 4
$ws.Range("Q35").Value = 0
$ws.Range("U35").Value = 0
$ws.Rows.Item(35).Hidden = $true

# --- Row 38 (becomes hidden) ---
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("P38").Value = 3
$ws.Range("Q38").Value = 0
$ws.Range("U38").Value = 0
$ws.Rows.Item(38).Hidden = $true

# --- Row 41 (becomes hidden) ---
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("P41").Value = 15
$ws.Range("Q41").Value = 0
$ws.Range("U41").Value = 0
$ws.Rows.Item(41).Hidden = $true

# --- Row 42 (stays visible; swaps talla data with row 44) ---
$ws.Range("C42").Value = "C20A60   "
$ws.Range("F42").Value = 2
$ws.Range("G42").Value = 45.48
$ws.Range("H42").Value = 18.19
$ws.Range("J42").Value = "AUMENTAR 25%"
$ws.Range("M42").Value = 90.95
$ws.Range("N42").Value = 54.57
$ws.Range("P42").Value = 1

# --- Row 43 (becomes hidden) ---
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("P43").Value = 7
$ws.Range("Q43").Value = 0
$ws.Range("U43").Value = 0
$ws.Rows.Item(43).Hidden = $true

# --- Row 44 (becomes hidden; swaps talla data with row 42) ---
$ws.Range("C44").Value = "C36A200  "
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 181.28
$ws.Range("H44").Value = 72.51
$ws.Range("J44").Value = "REDUCIR 50%"
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("P44").Value = 5
$ws.Range("Q44").Value = 0
$ws.Range("U44").Value = 0
$ws.Rows.Item(44).Hidden = $true

# --- Row 48 (stays visible) ---
$ws.Range("L48").Value = 0

# --- Totals ---
$ws.Range("C53").Value = 39
# Force text storage so the trailing currency symbol isn't re-parsed as a number
$ws.Range("C55").NumberFormat = "@"
$ws.Range("C55").Value = "921.04€"
$ws.Range("C64").Value = 0
